$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.80847
$ws.Range("H2").Value = 5.42541
$ws.Range("I2").Value = 0.1042562806587584
$ws.Range("J2").Value = 0.1042562806587584
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.94127933333334
$ws.Range("N2").Value = 263.823838
$ws.Range("O2").Value = 0.4109331243514438
$ws.Range("P2").Value = 0.4109331243514437
$ws.Range("Q2").Value = 159.0391654359534
$ws.Range("R2").Value = 1431.35248892358
$ws.Range("S2").Value = 0.04284235914436461
$ws.Range("T2").Value = 0.0428423591443646

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.80847
$ws.Range("H3").Value = 5.42541
$ws.Range("I3").Value = 0.1042562806587584
$ws.Range("J3").Value = 0.1042562806587584
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 52.441971
$ws.Range("N3").Value = 157.325913
$ws.Range("O3").Value = 0.2450515065683088
$ws.Range("P3").Value = 0.2450515065683087
$ws.Range("Q3").Value = 94.83973129437001
$ws.Range("R3").Value = 853.5575816493301
$ws.Range("S3").Value = 0.02554815864463719
$ws.Range("T3").Value = 0.02554815864463719

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.80847
$ws.Range("H4").Value = 5.42541
$ws.Range("I4").Value = 0.1042562806587584
$ws.Range("J4").Value = 0.1042562806587584
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.667459
$ws.Range("N4").Value = 164.002377
$ws.Range("O4").Value = 0.255450795093328
$ws.Range("P4").Value = 0.255450795093328
$ws.Range("Q4").Value = 98.86445957773
$ws.Range("R4").Value = 889.78013619957
$ws.Range("S4").Value = 0.026632349787753
$ws.Range("T4").Value = 0.02663234978775299

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.80847
$ws.Range("H5").Value = 5.42541
$ws.Range("I5").Value = 0.1042562806587584
$ws.Range("J5").Value = 0.1042562806587584
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 18.95316166666667
$ws.Range("N5").Value = 56.85948500000001
$ws.Range("O5").Value = 0.08856457398691947
$ws.Range("P5").Value = 0.08856457398691944
$ws.Range("Q5").Value = 34.27622427931667
$ws.Range("R5").Value = 308.4860185138501
$ws.Range("S5").Value = 0.009233413082003653
$ws.Range("T5").Value = 0.009233413082003651

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.255752333333334
$ws.Range("H6").Value = 18.767257
$ws.Range("I6").Value = 0.360637152397155
$ws.Range("J6").Value = 0.360637152397155
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 87.94127933333334
$ws.Range("N6").Value = 263.823838
$ws.Range("O6").Value = 0.4109331243514438
$ws.Range("P6").Value = 0.4109331243514437
$ws.Range("Q6").Value = 550.1388633858185
$ws.Range("R6").Value = 4951.249770472366
$ws.Range("S6").Value = 0.1481977517917707
$ws.Range("T6").Value = 0.1481977517917707

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.255752333333334
$ws.Range("H7").Value = 18.767257
$ws.Range("I7").Value = 0.360637152397155
$ws.Range("J7").Value = 0.360637152397155
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 52.441971
$ws.Range("N7").Value = 157.325913
$ws.Range("O7").Value = 0.2450515065683088
$ws.Range("P7").Value = 0.2450515065683087
$ws.Range("Q7").Value = 328.063982447849
$ws.Range("R7").Value = 2952.575842030642
$ws.Range("S7").Value = 0.08837467751942761
$ws.Range("T7").Value = 0.0883746775194276

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.255752333333334
$ws.Range("H8").Value = 18.767257
$ws.Range("I8").Value = 0.360637152397155
$ws.Range("J8").Value = 0.360637152397155
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.667459
$ws.Range("N8").Value = 164.002377
$ws.Range("O8").Value = 0.255450795093328
$ws.Range("P8").Value = 0.255450795093328
$ws.Range("Q8").Value = 341.9860841966544
$ws.Range("R8").Value = 3077.874757769889
$ws.Range("S8").Value = 0.09212504732004696
$ws.Range("T8").Value = 0.09212504732004695

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.255752333333334
$ws.Range("H9").Value = 18.767257
$ws.Range("I9").Value = 0.360637152397155
$ws.Range("J9").Value = 0.360637152397155
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.95316166666667
$ws.Range("N9").Value = 56.85948500000001
$ws.Range("O9").Value = 0.08856457398691947
$ws.Range("P9").Value = 0.08856457398691944
$ws.Range("Q9").Value = 118.5662853202939
$ws.Range("R9").Value = 1067.096567882645
$ws.Range("S9").Value = 0.03193967576590979
$ws.Range("T9").Value = 0.03193967576590978

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.749018666666667
$ws.Range("H10").Value = 11.247056
$ws.Range("I10").Value = 0.2161267493001954
$ws.Range("J10").Value = 0.2161267493001954
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 87.94127933333334
$ws.Range("N10").Value = 263.823838
$ws.Range("O10").Value = 0.4109331243514438
$ws.Range("P10").Value = 0.4109331243514437
$ws.Range("Q10").Value = 329.6934977912143
$ws.Range("R10").Value = 2967.241480120928
$ws.Range("S10").Value = 0.08881364034585051
$ws.Range("T10").Value = 0.0888136403458505

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.749018666666667
$ws.Range("H11").Value = 11.247056
$ws.Range("I11").Value = 0.2161267493001954
$ws.Range("J11").Value = 0.2161267493001954
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 52.441971
$ws.Range("N11").Value = 157.325913
$ws.Range("O11").Value = 0.2450515065683088
$ws.Range("P11").Value = 0.2450515065683087
$ws.Range("Q11").Value = 196.605928195792
$ws.Range("R11").Value = 1769.453353762128
$ws.Range("S11").Value = 0.05296218552572406
$ws.Range("T11").Value = 0.05296218552572404

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.749018666666667
$ws.Range("H12").Value = 11.247056
$ws.Range("I12").Value = 0.2161267493001954
$ws.Range("J12").Value = 0.2161267493001954
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.667459
$ws.Range("N12").Value = 164.002377
$ws.Range("O12").Value = 0.255450795093328
$ws.Range("P12").Value = 0.255450795093328
$ws.Range("Q12").Value = 204.9493242502347
$ws.Range("R12").Value = 1844.543918252112
$ws.Range("S12").Value = 0.05520974994967129
$ws.Range("T12").Value = 0.05520974994967127

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.749018666666667
$ws.Range("H13").Value = 11.247056
$ws.Range("I13").Value = 0.2161267493001954
$ws.Range("J13").Value = 0.2161267493001954
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 18.95316166666667
$ws.Range("N13").Value = 56.85948500000001
$ws.Range("O13").Value = 0.08856457398691947
$ws.Range("P13").Value = 0.08856457398691944
$ws.Range("Q13").Value = 71.05575688068446
$ws.Range("R13").Value = 639.5018119261601
$ws.Range("S13").Value = 0.01914117347894955
$ws.Range("T13").Value = 0.01914117347894954

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 5.533148
$ws.Range("H14").Value = 16.599444
$ws.Range("I14").Value = 0.3189798176438912
$ws.Range("J14").Value = 0.3189798176438912
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 87.94127933333334
$ws.Range("N14").Value = 263.823838
$ws.Range("O14").Value = 0.4109331243514438
$ws.Range("P14").Value = 0.4109331243514437
$ws.Range("Q14").Value = 486.5921138606747
$ws.Range("R14").Value = 4379.329024746072
$ws.Range("S14").Value = 0.131079373069458
$ws.Range("T14").Value = 0.131079373069458

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 5.533148
$ws.Range("H15").Value = 16.599444
$ws.Range("I15").Value = 0.3189798176438912
$ws.Range("J15").Value = 0.3189798176438912
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 52.441971
$ws.Range("N15").Value = 157.325913
$ws.Range("O15").Value = 0.2450515065683088
$ws.Range("P15").Value = 0.2450515065683087
$ws.Range("Q15").Value = 290.169186954708
$ws.Range("R15").Value = 2611.522682592372
$ws.Range("S15").Value = 0.07816648487851995
$ws.Range("T15").Value = 0.07816648487851992

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 5.533148
$ws.Range("H16").Value = 16.599444
$ws.Range("I16").Value = 0.3189798176438912
$ws.Range("J16").Value = 0.3189798176438912
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 54.667459
$ws.Range("N16").Value = 164.002377
$ws.Range("O16").Value = 0.255450795093328
$ws.Range("P16").Value = 0.255450795093328
$ws.Range("Q16").Value = 302.483141430932
$ws.Range("R16").Value = 2722.348272878387
$ws.Range("S16").Value = 0.0814836480358568
$ws.Range("T16").Value = 0.08148364803585677

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 5.533148
$ws.Range("H17").Value = 16.599444
$ws.Range("I17").Value = 0.3189798176438912
$ws.Range("J17").Value = 0.3189798176438912
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 18.95316166666667
$ws.Range("N17").Value = 56.85948500000001
$ws.Range("O17").Value = 0.08856457398691947
$ws.Range("P17").Value = 0.08856457398691944
$ws.Range("Q17").Value = 104.8706485695933
$ws.Range("R17").Value = 943.83583712634
$ws.Range("S17").Value = 0.02825031166005649
$ws.Range("T17").Value = 0.02825031166005647
